# Commit: "best combination and chains to investigate"
#
# 1. Rename the single existing sheet "Sheet 1" -> "All_Results"
# 2. Add a new sheet "Best_Combination" right after it, holding the same
#    7 headers (Chain, Start_Beta, Start_Alpha, Result_Alpha, Result_Beta,
#    Pr_Chi, Deviance) and the single best-result row that already exists
#    in All_Results (Chain 1 / 50 / 200 / 1849.94 / 186.2512 / ... ).

$wb = $excel.ActiveWorkbook

# --- 1. rename the original sheet -----------------------------------------
$allResults = $wb.Worksheets.Item(1)
$allResults.Name = "All_Results"

# --- 2. add the new sheet right after it -----------------------------------
$bestCombo = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $allResults)
$bestCombo.Name = "Best_Combination"

# --- header row --------------------------------------------------------
$headers = @("Chain", "Start_Beta", "Start_Alpha", "Result_Alpha", "Result_Beta", "Pr_Chi", "Deviance")
for ($col = 1; $col -le $headers.Length; $col++) {
    $bestCombo.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# --- best-combination data row ------------------------------------------
# same values as row 217 of All_Results (Start_Beta=50, Start_Alpha=200)
$bestCombo.Cells.Item(2, 1).Value = "Chain 1"
$bestCombo.Cells.Item(2, 2).Value = 50
$bestCombo.Cells.Item(2, 3).Value = 200

# Result_Alpha / Result_Beta are stored as text in the source data (they
# come out of the shared-string table rather than as numbers), so force
# text storage before assigning, then restore the default "Normal" style
# so no stray number-format style lingers on the cell.
$bestCombo.Cells.Item(2, 4).NumberFormat = "@"
$bestCombo.Cells.Item(2, 4).Value = "1849.94"
$bestCombo.Cells.Item(2, 4).Style = "Normal"

$bestCombo.Cells.Item(2, 5).NumberFormat = "@"
$bestCombo.Cells.Item(2, 5).Value = "186.2512"
$bestCombo.Cells.Item(2, 5).Style = "Normal"

$bestCombo.Cells.Item(2, 6).Value = 0.0000514336006743721
$bestCombo.Cells.Item(2, 7).Value = 19.7504377789202

# --- view cosmetics: keep both tabs at 100% zoom, All_Results active ------
$bestCombo.Activate()
$excel.ActiveWindow.Zoom = 100

$allResults.Activate()
$excel.ActiveWindow.Zoom = 100
